$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" '24.950.16'
Set-TextValue "E2" '  -0.12%  '
Set-TextValue "D3" '1.708.72'
Set-TextValue "E3" '  -0.45%  '
Set-TextValue "D4" '1.001'
Set-TextValue "E4" '  -0.18%  '
Set-TextValue "D5" '316.69'
Set-TextValue "E5" '  -0.37%  '
Set-TextValue "D6" '0.9999'
Set-TextValue "E6" '  -0.13%  '
Set-TextValue "D7" '0.4034'
Set-TextValue "E7" '  +1.53%  '
Set-TextValue "D8" '0.4081'
Set-TextValue "E8" '  -0.92%  '
Set-TextValue "B9" 'OKB'
Set-TextValue "C9" 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue "D9" '54.06'
Set-TextValue "E9" '  +1.11%  '
Set-TextValue "B10" 'Polygon'
Set-TextValue "C10" 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue "D10" '1.480'
Set-TextValue "E10" '  -3.13%  '
Set-TextValue "E11" '  -0.20%  '
Set-TextValue "D12" '0.08828'
Set-TextValue "E12" '  -1.33%  '
Set-TextValue "D13" '26.38'
Set-TextValue "E13" '  +5.91%  '
Set-TextValue "D14" '7.520'
Set-TextValue "E14" '  -2.51%  '
Set-TextValue "D15" '8.137'
Set-TextValue "E15" '  -0.30%  '
Set-TextValue "D16" '0.00001361'
Set-TextValue "E16" '  -0.84%  '
Set-TextValue "D17" '1.735.46'
Set-TextValue "E17" '  +3.16%  '
Set-TextValue "D18" '97.06'
Set-TextValue "E18" '  -3.46%  '
Set-TextValue "D19" '0.07159'
Set-TextValue "E19" '  +0.11%  '
Set-TextValue "D20" '21.14'
Set-TextValue "E20" '  +4.87%  '
Set-TextValue "D21" '7.278'
Set-TextValue "E21" '  -3.05%  '
Set-TextValue "D22" '1.002'
Set-TextValue "E22" '  -0.19%  '
Set-TextValue "D23" '14.40'
Set-TextValue "E23" '  -1.07%  '
Set-TextValue "D24" '24.952.88'
Set-TextValue "E24" '  -0.11%  '
Set-TextValue "D25" '2.320'
Set-TextValue "E25" '  -0.49%  '
Set-TextValue "D26" '2.912'
Set-TextValue "E26" '  -7.48%  '
Set-TextValue "D27" '23.35'
Set-TextValue "E27" '  +0.22%  '
Set-TextValue "D28" '6.266'
Set-TextValue "E28" '  +19.79%  '
Set-TextValue "D29" '167.14'
Set-TextValue "E29" '  +0.93%  '
Set-TextValue "D30" '146.96'
Set-TextValue "E30" '  +4.50%  '
Set-TextValue "D31" '8.402'
Set-TextValue "E31" '  -9.56%  '
Set-TextValue "D32" '1.920.38'
Set-TextValue "E32" '  +2.50%  '
Set-TextValue "D33" '2.223'
Set-TextValue "E33" '  +12.87%  '
Set-TextValue "D34" '0.08882'
Set-TextValue "E34" '  -1.69%  '
Set-TextValue "D35" '0.03224'
Set-TextValue "E35" '  +7.00%  '
Set-TextValue "D36" '7.284'
Set-TextValue "E36" '  -7.75%  '
Set-TextValue "D37" '1.030'
Set-TextValue "E37" '  -5.33%  '
Set-TextValue "D38" '0.2863'
Set-TextValue "E38" '  +1.88%  '
Set-TextValue "D39" '0.8483'
Set-TextValue "E39" '  +4.27%  '
Set-TextValue "D40" '10.87'
Set-TextValue "E40" '  -2.54%  '
Set-TextValue "D41" '0.09351'
Set-TextValue "E41" '  +0.25%  '
Set-TextValue "D42" '14.25'
Set-TextValue "E42" '  -2.47%  '
Set-TextValue "D43" '1.474'
Set-TextValue "E43" '  -1.03%  '
Set-TextValue "D44" '17.63'
Set-TextValue "E44" '  +4.97%  '
Set-TextValue "D45" '2.737'
Set-TextValue "E45" '  +3.30%  '
Set-TextValue "D46" '0.7445'
Set-TextValue "E46" '  +0.73%  '
Set-TextValue "D47" '4.245'
Set-TextValue "E47" '  -0.52%  '
Set-TextValue "D48" '1.399'
Set-TextValue "E48" '  +3.51%  '
Set-TextValue "E49" '  -0.14%  '
Set-TextValue "D50" '142.24'
Set-TextValue "E50" '  +0.90%  '
Set-TextValue "D51" '0.08394'
Set-TextValue "E51" '  +3.69%  '
